# R35 updates to about/resume
# Apply edits from later paragraphs to earlier ones so paragraph indices
# remain valid (the "About" paragraph edit inserts many new paragraphs
# and would shift every subsequent index).

$d = $word.ActiveDocument

# Paragraph 20: "Eclipse, Netbeans, GitHub, Brackets" bullet line - merge runs
$pBulletsEclipse = $d.Paragraphs.Item(20).Range
$xmlBulletsEclipse = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/>
        </w:rPr>
        <w:tab/>
        <w:t>Eclipse, Netbeans, GitHub, Brackets</w:t>
      </w:r>
    </w:p>
"@
$pBulletsEclipse.InsertXML($xmlBulletsEclipse)

# Paragraph 19: "Java , Javascript, jQuery, HTML, CSS, C++, MySQL, PHP" bullet line - merge runs
$pBulletsJava = $d.Paragraphs.Item(19).Range
$xmlBulletsJava = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:tab/>
        <w:t xml:space="preserve">Java </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:r>
        <w:t>Javascript, jQuery, HTML, CSS, C++, MySQL, PHP</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
    </w:p>
"@
$pBulletsJava.InsertXML($xmlBulletsJava)

# Paragraph 16: "Professional Skills" paragraph - merge the Java/HTML run and Eclipse/Netbeans run
$pProfSkills = $d.Paragraphs.Item(16).Range
$xmlProfSkills = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>Andrea has experience in Java and C++ as well as Javascript, jQuery, HTML, CSS, C++, MySQL, PHP</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">; employing a variety of tools to incude </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/>
        </w:rPr>
        <w:t>Eclipse, Netbeans, GitHub, Brackets</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Garamond" w:hAnsi="Garamond"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">She is familiar with Agile methodologies and holds the SAFe -Scaled Agile Framework Certificate.  </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">She is an excellent communicator who takes great pride in her work. </w:t>
      </w:r>
      <w:r>
        <w:t>She is adept at speaking to all types: technical professionals, business leaders, and clients</w:t>
      </w:r>
      <w:r>
        <w:t>. Additionally,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>s</w:t>
      </w:r>
      <w:r>
        <w:t>he speaks German, French, and Spanish with conversational proficiency</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
"@
$pProfSkills.InsertXML($xmlProfSkills)

# Paragraph 2: the underlined "About" heading paragraph gets replaced with the
# full new bio content (plain HTML-tagged text), ending in a plain "About" paragraph.
$pAbout = $d.Paragraphs.Item(2).Range
$xmlAbout = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>&lt;p&gt;</w:t>
      </w:r>
      <w:r>
        <w:t>I am</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> a software engineer new to the Nashville.  Having lived and worked in the Seattle area for the past six years, the warmth of the south </w:t>
      </w:r>
      <w:r>
        <w:t>called me</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> here to Nashville.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">                    &lt;br&gt;&lt;/br&gt;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">While in Seattle, </w:t>
      </w:r>
      <w:r>
        <w:t>I</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> worked as a Project Management Consultant for Sila Solutions Group.  During </w:t>
      </w:r>
      <w:r>
        <w:t>my</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> time at Sila, </w:t>
      </w:r>
      <w:r>
        <w:t>I</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> led many program level projects while consulting on the development of a large international airline customer service web-portal with a local original equipment manufacturer. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Additionally, I provided front-end development for a subsidiary of the original equipment manufacturer. </w:t>
      </w:r>
      <w:r>
        <w:t>I</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> also participated in the development of a full Linux/Apache/MySQL/PHP stack database: </w:t>
      </w:r>
      <w:r>
        <w:t>coding</w:t>
      </w:r>
      <w:r>
        <w:t>,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>gathering requirements,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> and</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> overseeing</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>general</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> progress</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> of the project</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">                        &lt;br&gt;&lt;/br&gt;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">Prior to </w:t>
      </w:r>
      <w:r>
        <w:t>my</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> time at Sila, </w:t>
      </w:r>
      <w:r>
        <w:t>I</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> taught Mathematics for a low income high school hoping to inspire under-privileged students to pursue challenging careers in math and computer science.  </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">                    &lt;/p&gt;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">                    &lt;p&gt;&lt;h3&gt;Education&lt;/h3&gt;&lt;/p&gt;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">                    &lt;p&gt;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">                    </w:t>
      </w:r>
      <w:r>
        <w:t>I hold</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> a Bachelors of Science in Computer Science and a Masters of Arts in Secondary Mathematics Instruction. </w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">                    &lt;/p&gt;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">                    &lt;p&gt;&lt;h3&gt;Professional Skills&lt;/h3&gt;&lt;/p&gt;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">                    &lt;p&gt;</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">                    </w:t>
      </w:r>
      <w:r>
        <w:t>I have</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> experience in Java, JavaS</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">cript, jQuery, HTML, CSS, C++, MySQL, PHP; employing a variety of tools to include Eclipse, Netbeans, GitHub, </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">and </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Brackets.  </w:t>
      </w:r>
      <w:r>
        <w:t>I am</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> familiar wi</w:t>
      </w:r>
      <w:r>
        <w:t>th Agile methodologies and hold</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> the SAFe -Scaled Agile Framework Certificate.  </w:t>
      </w:r>
      <w:r>
        <w:t>I am</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> an excellent communicator who takes great pride in </w:t>
      </w:r>
      <w:r>
        <w:t>my</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> work. </w:t>
      </w:r>
      <w:r>
        <w:t>I am</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> adept at speaking to all types: technical professionals, business leaders, and clients. Additionally, </w:t>
      </w:r>
      <w:r>
        <w:t>I speak</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> German, French, and Spanish with conversational proficiency.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
      </w:pPr>
      <w:r>
        <w:t>About</w:t>
      </w:r>
    </w:p>
    
"@
$pAbout.InsertXML($xmlAbout)

Write-Host "Edits applied"
